# Apply two changes described by the commit diff:
#  1. Slide 1: change the date text "January 2024" -> "May 2024"
#     (resulting in two runs "May" + " 2024", both inheriting the
#     original accent4 colored run formatting).
#  2. Slide 32: change the table's style id from
#     {69567550-1C8E-44CD-BECD-144A626A7D0F} to
#     {6ECCB492-A391-4F97-B320-6071CF3A2EC1}.

$p = $ppt.ActivePresentation

# --- Change 1: update the date on the title slide ---
$slide1 = $p.Slides.Item(1)
$subtitleShape = $slide1.Shapes.Item(2)
$dateParagraph = $subtitleShape.TextFrame.TextRange.Paragraphs(1, 1)

# Replace just the "January" portion (first 7 characters) with "May",
# leaving the trailing " 2024" (and its run formatting) untouched.
$monthRange = $dateParagraph.Characters(1, 7)
$monthRange.Text = "May"

# --- Change 2: update the table style on the resource table slide ---
$slide32 = $p.Slides.Item(32)
$tableShape = $slide32.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{6ECCB492-A391-4F97-B320-6071CF3A2EC1}")
